# Revert "Merge pull request #1 from FWSquatch/2score"
# This reverts the "2score" edits: the "Essential user" row should name
# jperalta (not brubble) as the user who was removed, and the
# "Users who should be added" row should name nscully (not bambam) as the
# user who was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Content revert -------------------------------------------------
# Row 4: "Essential user has been removed!" -> user is jperalta
$ws.Range("B4").Value = "jperalta"

# Row 6: "User nscully has been added." -> user is nscully
$ws.Range("B6").Value = "nscully"

# --- Cosmetic state restored by the revert --------------------------
# Selection moves from B12 back to A42.
[void]$ws.Range("A42").Select()

# Sheet-tab-bar split ratio back to ~0.5 (tabRatio 500/1000).
$excel.ActiveWindow.TabRatio = 0.5

# Column widths nudged back to their pre-"2score" values. Excel's
# ColumnWidth property is expressed in the default-font character units
# used on screen, which the file format stores with a constant ~5/6
# character offset baked in -- back that offset out so the persisted
# <col width="..."/> lands on the target value.
$offset = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 18.52 - $offset
$ws.Columns.Item(2).ColumnWidth = 50.33 - $offset
$ws.Columns.Item(3).ColumnWidth = 29.22 - $offset
$ws.Columns.Item(4).ColumnWidth = 23.53 - $offset
$ws.Columns.Item(5).ColumnWidth = 15.33 - $offset
$ws.Columns.Item(6).ColumnWidth = 54.64 - $offset
$ws.Columns.Item(7).ColumnWidth = 11.52 - $offset
